# Applies the "Added mean and variance values" commit to t580_project_data.xlsx
#  - A4:A6 become literal numbers (6) instead of the shared string `6"`
#  - C7:C15 becomes one shared-formula group (like the C17:C28 / C34:C41 groups)
#  - New summary block in columns G:I (rows 16-19) with per-sensor mean / St.Dev
#  - Two new header/label cells (H14/H15/I15) introducing "Mean" / "Variance" /
#    "Data between 3 sensors"
#  - New hidden `_xlchart.v1.*` defined names used by (now-untracked) charts
#  - Cosmetic: selection moved to H12

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. A4:A6 switch from the shared string `6"` to the literal number 6
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = 6
$ws.Range("A5").Value = 6
$ws.Range("A6").Value = 6

# ---------------------------------------------------------------------------
# 2. Re-key C7:C15 as a single shared formula (matches C17:C28 / C34:C41)
# ---------------------------------------------------------------------------
$ws.Range("C7:C15").Formula = "=((B7*0.001)/2)*13500"

# ---------------------------------------------------------------------------
# 3. New labels - written in the order that reproduces the target shared
#    string table (Mean=11, "Data between 3 sensors"=12, Variance=13)
# ---------------------------------------------------------------------------
$ws.Range("H15").Value = "Mean"
$ws.Range("H14").Value = "Data between 3 sensors"
$ws.Range("I15").Value = "Variance"

# ---------------------------------------------------------------------------
# 4. Per-sensor summary table: label (existing shared strings) + mean + stdev
# ---------------------------------------------------------------------------
$ws.Range("G16").Value = "6`""
$ws.Range("H16").Formula = "=AVERAGE(C4,C5,C6,C17,C18,C19,C30,C31,C32)"
$ws.Range("I16").Formula = "=STDEV.P(C4,C5,C6,C17,C19,C18,C30,C31,C32)"

$ws.Range("G17").Value = "12`""
$ws.Range("H17").Formula = "=AVERAGE(C7,C8,C9,C20,C21,C22,C33,C34,C35)"
$ws.Range("I17").Formula = "=STDEV.P(C7,C8,C9,C20,C21,C22,C33,C34,C35)"

$ws.Range("G18").Value = "18`""
$ws.Range("H18").Formula = "=AVERAGE(C38,C37,C36,C25,C24,C23,C12,C11,C10)"
$ws.Range("I18").Formula = "=STDEV.P(C10,C11,C12,C23,C24,C25,C36,C37,C38)"

$ws.Range("G19").Value = "24`""
$ws.Range("H19").Formula = "=AVERAGE(C13,C14,C15,C26,C27,C28,C39,C40,C41)"
$ws.Range("I19").Formula = "=STDEV.P(C15,C14,C13,C28,C27,C41,C40,C39,C26)"

# ---------------------------------------------------------------------------
# 5. Hidden defined names driving the (untracked) charts
# ---------------------------------------------------------------------------
$chartNames = @(
  @("_xlchart.v1.0", "(Sheet1!`$C`$10,Sheet1!`$C`$10:`$C`$12,Sheet1!`$C`$23:`$C`$25,Sheet1!`$C`$36:`$C`$38)"),
  @("_xlchart.v1.1", "(Sheet1!`$C`$13:`$C`$15,Sheet1!`$C`$26:`$C`$28,Sheet1!`$C`$39:`$C`$41)"),
  @("_xlchart.v1.2", "(Sheet1!`$C`$32,Sheet1!`$C`$30:`$C`$32,Sheet1!`$C`$17:`$C`$19,Sheet1!`$C`$4:`$C`$6)"),
  @("_xlchart.v1.3", "(Sheet1!`$C`$4,Sheet1!`$C`$4:`$C`$11,Sheet1!`$C`$12:`$C`$15,Sheet1!`$C`$17:`$C`$28,Sheet1!`$C`$30:`$C`$41)"),
  @("_xlchart.v1.4", "(Sheet1!`$C`$7:`$C`$9,Sheet1!`$C`$20:`$C`$22,Sheet1!`$C`$33:`$C`$35)"),
  @("_xlchart.v1.5", "(Sheet1!`$C`$10,Sheet1!`$C`$10:`$C`$12,Sheet1!`$C`$23:`$C`$25,Sheet1!`$C`$36:`$C`$38)"),
  @("_xlchart.v1.6", "(Sheet1!`$C`$13:`$C`$15,Sheet1!`$C`$26:`$C`$28,Sheet1!`$C`$39:`$C`$41)"),
  @("_xlchart.v1.7", "(Sheet1!`$C`$32,Sheet1!`$C`$30:`$C`$32,Sheet1!`$C`$17:`$C`$19,Sheet1!`$C`$4:`$C`$6)"),
  @("_xlchart.v1.8", "(Sheet1!`$C`$4,Sheet1!`$C`$4:`$C`$11,Sheet1!`$C`$12:`$C`$15,Sheet1!`$C`$17:`$C`$28,Sheet1!`$C`$30:`$C`$41)"),
  @("_xlchart.v1.9", "(Sheet1!`$C`$7:`$C`$9,Sheet1!`$C`$20:`$C`$22,Sheet1!`$C`$33:`$C`$35)")
)

foreach ($n in $chartNames) {
  $nm = $wb.Names.Add($n[0], "=" + $n[1])
  $nm.Visible = $false
}

# ---------------------------------------------------------------------------
# 6. Widen the two new columns (target raw widths 27 / 20.57 chars)
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 26.15
$ws.Columns.Item(9).ColumnWidth = 19.7

# ---------------------------------------------------------------------------
# 7. Cosmetic: zoom + selection move to H12 (matches saved view state)
# ---------------------------------------------------------------------------
$excel.ActiveWindow.Zoom = 85
$ws.Range("H12").Select()
